$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.316.89"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.281.22"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.54"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "264.96"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.76"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0930"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.87"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.53"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.620.75"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.857"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.280.34"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.198.53"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("E19").Value = "  -2.36%  "
$ws.Range("E20").Value = "  +2.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.45"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.53"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.14"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.64"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.34"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.92"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.34"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -8.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.35"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.25"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "172.06"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.35"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0909"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.78"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.13%  "
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("E37").Value = "  -2.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0354"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.86"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.90%  "
$ws.Range("E40").Value = "  -5.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.66"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +10.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "76.69"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.84"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.84%  "
$ws.Range("E44").Value = "  -5.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.20"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.38"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.65"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.60"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0992"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.25"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.87%  "
